# Apply the "corrected filenames for each algo" edit:
#  - sheet1 ("initial"): fill in the new V (col C) / T (col D) results for
#    buses 1-9, add an (empty) row 11, and touch up the selection.
#  - sheet2 ("line_imp"): no data changed upstream (only formatting noise).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("initial")

# ---- Column C (V) / Column D (T) values for rows 2-10 --------------------
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0

$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.14669763999999999

$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 0.045319579999999998

$ws1.Range("C5").Value = 1
$ws1.Range("D5").Value = -0.041805630000000003

$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = -0.084627740000000007

$ws1.Range("C7").Value = 1
$ws1.Range("D7").Value = -0.13192544

$ws1.Range("C8").Value = 1
$ws1.Range("D8").Value = 0.044822639999999997

$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = -0.01845047

$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = -0.0044904200000000002

# New trailing (blank) row, styled like the rest of column D
$ws1.Range("D11").Value = ""

# ---- Styling -----------------------------------------------------------
# Whole T column (D2:D11) is touched (applyNumberFormat), but only D3/D4
# switch to the 11pt font; D3 additionally gets justify/center alignment.
$ws1.Range("D2:D11").NumberFormat = "General"

$ws1.Range("D3").Font.Size = 11
$ws1.Range("D3").HorizontalAlignment = -4130
$ws1.Range("D3").VerticalAlignment = -4108

$ws1.Range("D4").Font.Size = 11

# ---- Sheet furniture -------------------------------------------------------
$ws1.Range("E20").Select()

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

Write-Host "edit applied"
